$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9377737641334534
$ws.Range("B1").Value = 2.009668588638306
$ws.Range("C1").Value = 2.998775243759155
$ws.Range("D1").Value = 3.638375997543335
$ws.Range("E1").Value = 1.779800534248352
